$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: was NCT03335345 / statut 1 -> becomes NCT04560153 / statut 4 (Karate study) ---
$ws.Range("C2").Value = "NCT04560153"
$ws.Range("G2").Value = "Interest of Karate Kata Practice on the Self-esteem of Patients Living With HIV"
$ws.Range("H2").Value = "KATACHRO"
# F2 stays "2020", I2 stays "OTHER"

# --- Row 3: was NCT04560153 / statut 4 -> becomes NCT03335345 / statut 1 (Non-inferiority study) ---
$ws.Range("C3").Value = "NCT03335345"
$ws.Range("G3").Value = "Non-inferiority Study of the Pursuit of Enteral Nutrition Compared to a Strategy of Gastric Emptiness Peri-extubation. Cluster Randomized Trial"
$ws.Range("H3").Value = "AMBROISIE"
# F3 stays "2020", I3 stays "OTHER"

# --- Row 4: unchanged content (NCT03071601 / statut 1), only statut_name wording changes below ---

# --- Row 5 (new): NCT05627167 / statut 4 (DC-SCENIC) ---
$ws.Range("C5").Value = "NCT05627167"
$ws.Range("G5").Value = "Daytime Cyclic Enteral Nutrition Versus Standard Continuous Enteral Nutrition in the Intensive Care Unit: a Pilot Randomized Controlled Trial"
$ws.Range("H5").Value = "DC-SCENIC"
$ws.Range("I5").Value = "OTHER"

# --- statut (column A) + statut_name (column B) for all 4 data rows, and the new
# completion_year cell (F5) --- these must stay text values ("1"/"4"/"2025"), not
# auto-convert to numbers, so the cells are pre-formatted as Text before the value
# is entered (matches how the sibling cells in the same columns are already stored).
$ws.Range("A2:A5").NumberFormat = "@"
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "2025"

$ws.Range("A2").Value = "4"
$ws.Range("B2").Value = "4 : pas de résultats postés ni publiés"

$ws.Range("A3").Value = "1"
$ws.Range("B3").Value = "1 : résultats postés ou publiés dans les 12 mois"

$ws.Range("A4").Value = "1"
$ws.Range("B4").Value = "1 : résultats postés ou publiés dans les 12 mois"

$ws.Range("A5").Value = "4"
$ws.Range("B5").Value = "4 : pas de résultats postés ni publiés"
